# Add data for 2025-11-01
# Updates the 2025 (column L) crime-count figures across the citywide
# summary, the by-neighborhood summary, and every individual neighborhood
# sheet affected by the new day of data. A few columns for 2022 (I) and
# 2023 (J) also received small corrections in some sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 5614
$ws.Range('L3').Value = 6097
$ws.Range('I4').Value = 1849
$ws.Range('J4').Value = 1878
$ws.Range('L4').Value = 1501
$ws.Range('L5').Value = 363
$ws.Range('L6').Value = 5003
$ws.Range('I7').Value = 26320
$ws.Range('J7').Value = 29355
$ws.Range('L7').Value = 18578

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 362
$ws.Range('L3').Value = 436
$ws.Range('L6').Value = 306
$ws.Range('L7').Value = 1234

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 232
$ws.Range('L6').Value = 240
$ws.Range('L7').Value = 849

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L3').Value = 248
$ws.Range('L7').Value = 711

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('L3').Value = 132
$ws.Range('L6').Value = 72
$ws.Range('L7').Value = 325

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('L6').Value = 35
$ws.Range('L7').Value = 83

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L7').Value = 599
$ws.Range('L8').Value = 1234
$ws.Range('L10').Value = 124
$ws.Range('I11').Value = 390
$ws.Range('L11').Value = 305
$ws.Range('L14').Value = 95
$ws.Range('L15').Value = 152
$ws.Range('L19').Value = 509
$ws.Range('L20').Value = 458
$ws.Range('L23').Value = 205
$ws.Range('L27').Value = 163
$ws.Range('J29').Value = 1552
$ws.Range('L29').Value = 1049
$ws.Range('L30').Value = 83
$ws.Range('L33').Value = 849
$ws.Range('L37').Value = 711
$ws.Range('L41').Value = 78
$ws.Range('L42').Value = 604
$ws.Range('L43').Value = 135
$ws.Range('L44').Value = 126
$ws.Range('L45').Value = 35
$ws.Range('L47').Value = 122
$ws.Range('L50').Value = 92
$ws.Range('L51').Value = 231
$ws.Range('L54').Value = 409
$ws.Range('L55').Value = 195
$ws.Range('L63').Value = 54
$ws.Range('L66').Value = 52
$ws.Range('L67').Value = 642
$ws.Range('L68').Value = 58
$ws.Range('L79').Value = 505
$ws.Range('L81').Value = 16
$ws.Range('L85').Value = 922
$ws.Range('L86').Value = 126
$ws.Range('L88').Value = 197
$ws.Range('L91').Value = 246
$ws.Range('L94').Value = 230
$ws.Range('L96').Value = 211
$ws.Range('L99').Value = 325
$ws.Range('I101').Value = 26320
$ws.Range('J101').Value = 29355
$ws.Range('L101').Value = 18578

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L2').Value = 189
$ws.Range('L7').Value = 642

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L2').Value = 74
$ws.Range('L6').Value = 196
$ws.Range('L7').Value = 409

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L2').Value = 308
$ws.Range('L3').Value = 403
$ws.Range('J4').Value = 84
$ws.Range('L4').Value = 57
$ws.Range('L6').Value = 264
$ws.Range('J7').Value = 1552
$ws.Range('L7').Value = 1049

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L2').Value = 180
$ws.Range('L3').Value = 157
$ws.Range('L7').Value = 509

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('L3').Value = 36
$ws.Range('L7').Value = 126

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('L6').Value = 24
$ws.Range('L7').Value = 95

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('L3').Value = 28
$ws.Range('L7').Value = 78

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L2').Value = 164
$ws.Range('L6').Value = 165
$ws.Range('L7').Value = 604

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('L5').Value = 2
$ws.Range('L7').Value = 124

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('L2').Value = 58
$ws.Range('L7').Value = 195

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('L3').Value = 82
$ws.Range('L7').Value = 205

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L2').Value = 66
$ws.Range('L7').Value = 211

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('L2').Value = 84
$ws.Range('L7').Value = 246

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L2').Value = 165
$ws.Range('L6').Value = 128
$ws.Range('L7').Value = 505

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L2').Value = 142
$ws.Range('L3').Value = 156
$ws.Range('L7').Value = 458

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L2').Value = 206
$ws.Range('L7').Value = 599

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('L4').Value = 30
$ws.Range('L6').Value = 88
$ws.Range('L7').Value = 230

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('L2').Value = 47
$ws.Range('L7').Value = 122

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('L3').Value = 48
$ws.Range('L7').Value = 152

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('L4').Value = 13
$ws.Range('L7').Value = 92

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('L3').Value = 16
$ws.Range('L7').Value = 52

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('I4').Value = 35
$ws.Range('L4').Value = 23
$ws.Range('I7').Value = 390
$ws.Range('L7').Value = 305

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('L2').Value = 60
$ws.Range('L3').Value = 69
$ws.Range('L7').Value = 197

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('L4').Value = 21
$ws.Range('L6').Value = 49
$ws.Range('L7').Value = 163

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('L4').Value = 67
$ws.Range('L7').Value = 126

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('L5').Value = 5
$ws.Range('L7').Value = 231

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('L6').Value = 16
$ws.Range('L7').Value = 58

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('L3').Value = 43
$ws.Range('L7').Value = 135

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L2').Value = 280
$ws.Range('L6').Value = 190
$ws.Range('L7').Value = 922

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range('L3').Value = 15
$ws.Range('L7').Value = 35

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range('L2').Value = 7
$ws.Range('L7').Value = 16
